$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.921.11'
$ws.Range("E2").Value = '  -3.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.737.46'
$ws.Range("E3").Value = '  -1.13%  '

$ws.Range("E4").Value = '  -0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.34'
$ws.Range("E5").Value = '  -5.60%  '

$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4981'
$ws.Range("E7").Value = '  +3.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3543'
$ws.Range("E8").Value = '  +0.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.43'
$ws.Range("E9").Value = '  -2.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07237'
$ws.Range("E10").Value = '  -4.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.058'
$ws.Range("E11").Value = '  -1.88%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.90'
$ws.Range("E13").Value = '  -2.83%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.951'
$ws.Range("E14").Value = '  -2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.733.15'
$ws.Range("E15").Value = '  -1.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.836'
$ws.Range("E16").Value = '  -4.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.32'
$ws.Range("E17").Value = '  -6.43%  '

$ws.Range("E18").Value = '  -5.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06386'
$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("E21").Value = '  -1.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.730'
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '26.986.00'
$ws.Range("E23").Value = '  -2.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.19'
$ws.Range("E24").Value = '  +0.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.046'
$ws.Range("E25").Value = '  -5.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.22'
$ws.Range("E26").Value = '  -6.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.83'
$ws.Range("E27").Value = '  -0.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.934.18'
$ws.Range("E28").Value = '  -1.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.123'
$ws.Range("E29").Value = '  -2.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.58'
$ws.Range("E30").Value = '  -1.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.056'
$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09420'
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.573'
$ws.Range("E33").Value = '  -2.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.366'
$ws.Range("E34").Value = '  -3.06%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02188'
$ws.Range("E35").Value = '  -3.12%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05902'
$ws.Range("E36").Value = '  -1.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.03'
$ws.Range("E37").Value = '  -4.75%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2001'
$ws.Range("E38").Value = '  -2.88%  '

$ws.Range("B39").Value = 'WEMIXTOKEN'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.424'
$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.753'
$ws.Range("E40").Value = '  -2.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5993'
$ws.Range("E42").Value = '  -2.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.103'
$ws.Range("E43").Value = '  -6.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.439'
$ws.Range("E44").Value = '  -3.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.86'
$ws.Range("E45").Value = '  -2.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.572'
$ws.Range("E46").Value = '  -4.41%  '

$ws.Range("E47").Value = '  -2.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.77'
$ws.Range("E48").Value = '  -2.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.850'
$ws.Range("E49").Value = '  -3.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06657'
$ws.Range("E50").Value = '  -1.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.097'
$ws.Range("E51").Value = '  -4.25%  '
